$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Combine the old rows 2-5 into a single Python-tuple-like string in A2
$ws.Range("A2").Value = "('Memory Lapse', ['{1}{U}', 'Instant', 'Counter target spell. If that spell is countered this way, put it on top of its owner" + [char]8217 + "s library instead of into that player" + [char]8217 + "s graveyard.'])"

# Remove old rows 3-5 which are no longer needed
$ws.Range("A3:A5").ClearContents()
